$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030850055578532
$ws.Range("D2").Value = 1.039857862803707
$ws.Range("E2").Value = 1.040413663790888
$ws.Range("F2").Value = 1.051178108451574
$ws.Range("I2").Value = 1.033239771132998
$ws.Range("J2").Value = 1.035989025708236
$ws.Range("K2").Value = 1.042641858123046
$ws.Range("L2").Value = 1.043196081642828
$ws.Range("M2").Value = 1.053930331545584
$ws.Range("N2").Value = 1.015975833862661

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031720252409871
$ws.Range("D3").Value = 1.040530652899418
$ws.Range("E3").Value = 1.041182124979617
$ws.Range("F3").Value = 1.052018592578164
$ws.Range("I3").Value = 1.033354785234799
$ws.Range("J3").Value = 1.036501254296195
$ws.Range("K3").Value = 1.043125315135837
$ws.Range("L3").Value = 1.043775071598247
$ws.Range("M3").Value = 1.054583338220683
$ws.Range("N3").Value = 1.016147445345857

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032284044796518
$ws.Range("D4").Value = 1.040966413693967
$ws.Range("E4").Value = 1.041680369225398
$ws.Range("F4").Value = 1.052563452405082
$ws.Range("I4").Value = 1.033427986529348
$ws.Range("J4").Value = 1.036832775102355
$ws.Range("K4").Value = 1.043437883382169
$ws.Range("L4").Value = 1.04415005322024
$ws.Range("M4").Value = 1.055006239063886
$ws.Range("N4").Value = 1.016258462200681

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032521233165976
$ws.Range("D5").Value = 1.041149706380773
$ws.Range("E5").Value = 1.041890068310819
$ws.Range("F5").Value = 1.052792751399173
$ws.Range("I5").Value = 1.033458467479665
$ws.Range("J5").Value = 1.036972163035423
$ws.Range("K5").Value = 1.043569223233221
$ws.Range("L5").Value = 1.044307774580904
$ws.Range("M5").Value = 1.055184111549561
$ws.Range("N5").Value = 1.016305126682336

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032561068063292
$ws.Range("D6").Value = 1.041180487772208
$ws.Range("E6").Value = 1.041925291566497
$ws.Range("F6").Value = 1.052831265721213
$ws.Range("I6").Value = 1.033463568171356
$ws.Range("J6").Value = 1.036995567817923
$ws.Range("K6").Value = 1.043591271987507
$ws.Range("L6").Value = 1.044334261286912
$ws.Range("M6").Value = 1.055213982057104
$ws.Range("N6").Value = 1.016312961429442

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032287213452894
$ws.Range("D7").Value = 1.040968862473737
$ws.Range("E7").Value = 1.041683170305093
$ws.Range("F7").Value = 1.052566515369205
$ws.Range("I7").Value = 1.03342839496894
$ws.Range("J7").Value = 1.036834637547183
$ws.Range("K7").Value = 1.043439638604308
$ws.Range("L7").Value = 1.044152160391845
$ws.Range("M7").Value = 1.055008615472952
$ws.Range("N7").Value = 1.016259085761894

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031143993301843
$ws.Range("D8").Value = 1.040085147165547
$ws.Range("E8").Value = 1.040673161006613
$ws.Range("F8").Value = 1.051461943605221
$ws.Range("I8").Value = 1.033278892995551
$ws.Range("J8").Value = 1.036162119525475
$ws.Range("K8").Value = 1.042805297829701
$ws.Range("L8").Value = 1.043391683483018
$ws.Range("M8").Value = 1.054150942286393
$ws.Range("N8").Value = 1.016033836120499

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029135052822174
$ws.Range("D9").Value = 1.038531230908694
$ws.Range("E9").Value = 1.03890112364454
$ws.Range("F9").Value = 1.049523372469711
$ws.Range("I9").Value = 1.033006133728351
$ws.Range("J9").Value = 1.034977685545998
$ws.Range("K9").Value = 1.041685573535219
$ws.Range("L9").Value = 1.042054266536998
$ws.Range("M9").Value = 1.052642457290749
$ws.Range("N9").Value = 1.015636729995198

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027799590175084
$ws.Range("D10").Value = 1.037497617633991
$ws.Range("E10").Value = 1.037725070565595
$ws.Range("F10").Value = 1.048236370478044
$ws.Range("I10").Value = 1.032818070021963
$ws.Range("J10").Value = 1.034188562270421
$ws.Range("K10").Value = 1.040937873916829
$ws.Range("L10").Value = 1.041164516967543
$ws.Range("M10").Value = 1.051638812162824
$ws.Range("N10").Value = 1.015371895135244

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02722224830699
$ws.Range("D11").Value = 1.037050627793609
$ws.Range("E11").Value = 1.037217107989884
$ws.Range("F11").Value = 1.047680384964509
$ws.Range("I11").Value = 1.032735168586319
$ws.Range("J11").Value = 1.033846997778998
$ws.Range("K11").Value = 1.040613840861686
$ws.Range("L11").Value = 1.040779705379985
$ws.Range("M11").Value = 1.051204720718871
$ws.Range("N11").Value = 1.01525720163446

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027007937953928
$ws.Range("D12").Value = 1.036884683611713
$ws.Range("E12").Value = 1.037028621673934
$ws.Range("F12").Value = 1.047474063732747
$ws.Range("I12").Value = 1.032704155329147
$ws.Range("J12").Value = 1.033720146511567
$ws.Range("K12").Value = 1.040493440787289
$ws.Range("L12").Value = 1.04063683917815
$ws.Range("M12").Value = 1.051043555657891
$ws.Range("N12").Value = 1.015214597102883

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027053901886175
$ws.Range("D13").Value = 1.03692027520112
$ws.Range("E13").Value = 1.037069043837614
$ws.Range("F13").Value = 1.04751831141111
$ws.Range("I13").Value = 1.032710817721566
$ws.Range("J13").Value = 1.033747355568236
$ws.Range("K13").Value = 1.040519268788792
$ws.Range("L13").Value = 1.04066748127419
$ws.Range("M13").Value = 1.051078122652103
$ws.Range("N13").Value = 1.015223736013965

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027204530478837
$ws.Range("D14").Value = 1.037036908994882
$ws.Range("E14").Value = 1.037201523691972
$ws.Range("F14").Value = 1.047663326366097
$ws.Range("I14").Value = 1.032732609504311
$ws.Range("J14").Value = 1.033836511777909
$ws.Range("K14").Value = 1.040603889357047
$ws.Range("L14").Value = 1.040767894579604
$ws.Range("M14").Value = 1.051191397208027
$ws.Range("N14").Value = 1.015253679973865

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027297356340588
$ws.Range("D15").Value = 1.037108782608188
$ws.Range("E15").Value = 1.037283174612699
$ws.Range("F15").Value = 1.047752700969891
$ws.Range("I15").Value = 1.032746007012384
$ws.Range("J15").Value = 1.033891446681779
$ws.Range("K15").Value = 1.040656021654611
$ws.Range("L15").Value = 1.040829771844862
$ws.Range("M15").Value = 1.051261199498918
$ws.Range("N15").Value = 1.015272129151026

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027837926017876
$ws.Range("D16").Value = 1.037527295079665
$ws.Range("E16").Value = 1.037758809416023
$ws.Range("F16").Value = 1.048273296857723
$ws.Range("I16").Value = 1.032823541032402
$ws.Range("J16").Value = 1.034211233653553
$ws.Range("K16").Value = 1.040959373268233
$ws.Range("L16").Value = 1.041190065374902
$ws.Range("M16").Value = 1.05166763192392
$ws.Range("N16").Value = 1.015379506621833

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028177259114381
$ws.Range("D17").Value = 1.037789971239318
$ws.Range("E17").Value = 1.03805750555037
$ws.Range("F17").Value = 1.048600200987804
$ws.Range("I17").Value = 1.032871783358359
$ws.Range("J17").Value = 1.034411863585831
$ws.Range("K17").Value = 1.041149585318473
$ws.Range("L17").Value = 1.041416191068389
$ws.Range("M17").Value = 1.051922709564143
$ws.Range("N17").Value = 1.015446857126168

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028375275249292
$ws.Range("D18").Value = 1.037943240817268
$ws.Range("E18").Value = 1.038231852958268
$ws.Range("F18").Value = 1.048791003434939
$ws.Range("I18").Value = 1.032899780546418
$ws.Range("J18").Value = 1.034528900192628
$ws.Range("K18").Value = 1.041260506238819
$ws.Range("L18").Value = 1.041548130243797
$ws.Range("M18").Value = 1.052071539540712
$ws.Range("N18").Value = 1.015486139741689

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028442808662358
$ws.Range("D19").Value = 1.037995511057271
$ws.Range("E19").Value = 1.038291321725851
$ws.Range("F19").Value = 1.048856083245587
$ws.Range("I19").Value = 1.032909302801367
$ws.Range("J19").Value = 1.034568808742259
$ws.Range("K19").Value = 1.041298322859485
$ws.Range("L19").Value = 1.041593125491792
$ws.Range("M19").Value = 1.052122294722991
$ws.Range("N19").Value = 1.015499533779357

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028140842679909
$ws.Range("D20").Value = 1.037761782868965
$ws.Range("E20").Value = 1.038025445545363
$ws.Range("F20").Value = 1.048565114345416
$ws.Range("I20").Value = 1.032866622069812
$ws.Range("J20").Value = 1.034390336584162
$ws.Range("K20").Value = 1.041129180101592
$ws.Range("L20").Value = 1.041391925371761
$ws.Range("M20").Value = 1.051895337221604
$ws.Range("N20").Value = 1.015439631236535

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027160170231974
$ws.Range("D21").Value = 1.037002560812295
$ws.Range("E21").Value = 1.037162506317457
$ws.Range("F21").Value = 1.04762061763458
$ws.Range("I21").Value = 1.032726198434385
$ws.Range("J21").Value = 1.033810256906396
$ws.Range("K21").Value = 1.040578971795232
$ws.Range("L21").Value = 1.040738323420982
$ws.Range("M21").Value = 1.051158038567827
$ws.Range("N21").Value = 1.015244862282469

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026544394025279
$ws.Range("D22").Value = 1.036525716378717
$ws.Range("E22").Value = 1.03662106322021
$ws.Range("F22").Value = 1.047027913197297
$ws.Range("I22").Value = 1.032636636158596
$ws.Range("J22").Value = 1.033445659890074
$ws.Range("K22").Value = 1.040232804601601
$ws.Range("L22").Value = 1.040327783401285
$ws.Range("M22").Value = 1.050694909198325
$ws.Range("N22").Value = 1.015122390482339

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026870751201319
$ws.Range("D23").Value = 1.036778451797512
$ws.Range("E23").Value = 1.036907985563699
$ws.Range("F23").Value = 1.047342008481855
$ws.Range("I23").Value = 1.032684235197164
$ws.Range("J23").Value = 1.03363892765989
$ws.Range("K23").Value = 1.040416335649441
$ws.Range("L23").Value = 1.040545379509131
$ws.Range("M23").Value = 1.050940380612848
$ws.Range("N23").Value = 1.01518731615513

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02815729742467
$ws.Range("D24").Value = 1.037774519810412
$ws.Range("E24").Value = 1.038039931701787
$ws.Range("F24").Value = 1.048580968103468
$ws.Range("I24").Value = 1.032868954671987
$ws.Range("J24").Value = 1.03440006367094
$ws.Range("K24").Value = 1.041138400422976
$ws.Range("L24").Value = 1.041402889860148
$ws.Range("M24").Value = 1.051907705458394
$ws.Range("N24").Value = 1.015442896310629

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029653743681105
$ws.Range("D25").Value = 1.038932552726979
$ws.Range("E25").Value = 1.039358310701113
$ws.Range("F25").Value = 1.050023600123247
$ws.Range("I25").Value = 1.033077748561849
$ws.Range("J25").Value = 1.035283807715871
$ws.Range("K25").Value = 1.041975269721681
$ws.Range("L25").Value = 1.042399699056034
$ws.Range("M25").Value = 1.053032089512081
$ws.Range("N25").Value = 1.015739410715238
